# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record is inserted at row 419 of the "Ciboulette" sheet,
# pushing all subsequent rows (419-515) down by one (to 420-516).
# The new row 419 keeps the same static attributes as the (old) row 419 but
# carries a new date (D) and a new minimum-volume figure (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 419 - this shifts rows 419:515 down to 420:516
# and carries formatting (e.g. the date-number style on column D) along.
$ws.Rows.Item(419).Insert()

# Populate the newly-inserted (now blank) row 419 with its values.
$ws.Range("A419").Value = 9
$ws.Range("B419").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C419").Value = "Metropolitana"
$ws.Range("D419").Value = 44943
$ws.Range("E419").Value = 13
$ws.Range("F419").Value = 100112039
$ws.Range("G419").Value = "Ciboulette"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 340
$ws.Range("K419").Value = 1000
$ws.Range("L419").Value = 1000
$ws.Range("M419").Value = 1000
$ws.Range("N419").Value = "$/docena de atados"
$ws.Range("O419").Value = "Región Metropolitana"
$ws.Range("P419").Value = 333
$ws.Range("Q419").Value = 3
$ws.Range("R419").Value = "Hortaliza"
